# Insert one new data row before the current row 99 (shifts old rows 99-157 down to 100-158)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(99).Insert()

# Populate the newly inserted row 99 with the new record (same template columns as every
# other row in this sheet: A/B/C/E/F/G/H/I/R), plus the record-specific values.
$ws.Range("A99").Value = 4
$ws.Range("B99").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C99").Value = "Los Lagos"
$ws.Range("D99").Value = 44438
$ws.Range("E99").Value = 10
$ws.Range("F99").Value = 100112040
$ws.Range("G99").Value = "Cilantro"
$ws.Range("H99").Value = "Sin especificar"
$ws.Range("I99").Value = "Primera"
$ws.Range("J99").Value = 150
$ws.Range("K99").Value = 15000
$ws.Range("L99").Value = 15000
$ws.Range("M99").Value = 15000
$ws.Range("N99").Value = "$/caja 36 atados"
$ws.Range("O99").Value = "Región Metropolitana"
$ws.Range("P99").Value = 417
$ws.Range("Q99").Value = 36
$ws.Range("R99").Value = "Hortaliza"
